$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 7).Value = 3.027114666666666
$ws.Cells.Item(2, 8).Value = 9.081344
$ws.Cells.Item(2, 9).Value = 0.207506525262911
$ws.Cells.Item(2, 10).Value = 0.207506525262911
$ws.Cells.Item(2, 13).Value = 18.95273633333333
$ws.Cells.Item(2, 14).Value = 56.858209
$ws.Cells.Item(2, 15).Value = 0.08721078561875104
$ws.Cells.Item(2, 16).Value = 0.08721078561875105
$ws.Cells.Item(2, 17).Value = 57.37210612809954
$ws.Cells.Item(2, 18).Value = 516.3489551528959
$ws.Cells.Item(2, 19).Value = 0.01809680708919568
$ws.Cells.Item(2, 20).Value = 0.01809680708919568

$ws.Cells.Item(3, 7).Value = 3.027114666666666
$ws.Cells.Item(3, 8).Value = 9.081344
$ws.Cells.Item(3, 9).Value = 0.207506525262911
$ws.Cells.Item(3, 10).Value = 0.207506525262911
$ws.Cells.Item(3, 15).Value = 0.04852204497892696
$ws.Cells.Item(3, 16).Value = 0.04852204497892696
$ws.Cells.Item(3, 17).Value = 31.92050036394666
$ws.Cells.Item(3, 18).Value = 287.28450327552
$ws.Cells.Item(3, 19).Value = 0.01006864095222781
$ws.Cells.Item(3, 20).Value = 0.01006864095222781

$ws.Cells.Item(4, 7).Value = 3.027114666666666
$ws.Cells.Item(4, 8).Value = 9.081344
$ws.Cells.Item(4, 9).Value = 0.207506525262911
$ws.Cells.Item(4, 10).Value = 0.207506525262911
$ws.Cells.Item(4, 13).Value = 101.4555613333333
$ws.Cells.Item(4, 14).Value = 304.366684
$ws.Cells.Item(4, 15).Value = 0.4668465309523581
$ws.Cells.Item(4, 16).Value = 0.4668465309523581
$ws.Cells.Item(4, 17).Value = 307.1176177270329
$ws.Cells.Item(4, 18).Value = 2764.058559543296
$ws.Cells.Item(4, 19).Value = 0.09687370146896787
$ws.Cells.Item(4, 20).Value = 0.09687370146896787

$ws.Cells.Item(5, 7).Value = 3.027114666666666
$ws.Cells.Item(5, 8).Value = 9.081344
$ws.Cells.Item(5, 9).Value = 0.207506525262911
$ws.Cells.Item(5, 10).Value = 0.207506525262911
$ws.Cells.Item(5, 13).Value = 2.410466333333333
$ws.Cells.Item(5, 14).Value = 7.231399000000001
$ws.Cells.Item(5, 15).Value = 0.01109173150200089
$ws.Cells.Item(5, 16).Value = 0.01109173150200089
$ws.Cells.Item(5, 17).Value = 7.296757991139555
$ws.Cells.Item(5, 18).Value = 65.670821920256
$ws.Cells.Item(5, 19).Value = 0.002301606663129375
$ws.Cells.Item(5, 20).Value = 0.002301606663129375

$ws.Cells.Item(6, 7).Value = 3.027114666666666
$ws.Cells.Item(6, 8).Value = 9.081344
$ws.Cells.Item(6, 9).Value = 0.207506525262911
$ws.Cells.Item(6, 10).Value = 0.207506525262911
$ws.Cells.Item(6, 13).Value = 83.95738966666666
$ws.Cells.Item(6, 14).Value = 251.872169
$ws.Cells.Item(6, 15).Value = 0.386328906947963
$ws.Cells.Item(6, 16).Value = 0.386328906947963
$ws.Cells.Item(6, 17).Value = 254.1486456350151
$ws.Cells.Item(6, 18).Value = 2287.337810715136
$ws.Cells.Item(6, 19).Value = 0.08016576908939028
$ws.Cells.Item(6, 20).Value = 0.08016576908939028

$ws.Cells.Item(7, 8).Value = 7.555711000000001
$ws.Cells.Item(7, 9).Value = 0.1726461783080517
$ws.Cells.Item(7, 10).Value = 0.1726461783080516
$ws.Cells.Item(7, 13).Value = 18.95273633333333
$ws.Cells.Item(7, 14).Value = 56.858209
$ws.Cells.Item(7, 15).Value = 0.08721078561875104
$ws.Cells.Item(7, 16).Value = 0.08721078561875105
$ws.Cells.Item(7, 17).Value = 47.7337994646221
$ws.Cells.Item(7, 18).Value = 429.604195181599
$ws.Cells.Item(7, 19).Value = 0.01505660884432016
$ws.Cells.Item(7, 20).Value = 0.01505660884432016

$ws.Cells.Item(8, 8).Value = 7.555711000000001
$ws.Cells.Item(8, 9).Value = 0.1726461783080517
$ws.Cells.Item(8, 10).Value = 0.1726461783080516
$ws.Cells.Item(8, 15).Value = 0.04852204497892696
$ws.Cells.Item(8, 16).Value = 0.04852204497892696
$ws.Cells.Item(8, 19).Value = 0.008377145629303126
$ws.Cells.Item(8, 20).Value = 0.008377145629303124

$ws.Cells.Item(9, 8).Value = 7.555711000000001
$ws.Cells.Item(9, 9).Value = 0.1726461783080517
$ws.Cells.Item(9, 10).Value = 0.1726461783080516
$ws.Cells.Item(9, 13).Value = 101.4555613333333
$ws.Cells.Item(9, 14).Value = 304.366684
$ws.Cells.Item(9, 15).Value = 0.4668465309523581
$ws.Cells.Item(9, 16).Value = 0.4668465309523581
$ws.Cells.Item(9, 17).Value = 255.5229669258138
$ws.Cells.Item(9, 18).Value = 2299.706702332324
$ws.Cells.Item(9, 19).Value = 0.08059926942529617
$ws.Cells.Item(9, 20).Value = 0.08059926942529616

$ws.Cells.Item(10, 8).Value = 7.555711000000001
$ws.Cells.Item(10, 9).Value = 0.1726461783080517
$ws.Cells.Item(10, 10).Value = 0.1726461783080516
$ws.Cells.Item(10, 13).Value = 2.410466333333333
$ws.Cells.Item(10, 14).Value = 7.231399000000001
$ws.Cells.Item(10, 15).Value = 0.01109173150200089
$ws.Cells.Item(10, 16).Value = 0.01109173150200089
$ws.Cells.Item(10, 17).Value = 6.070928996632111
$ws.Cells.Item(10, 18).Value = 54.63836096968901
$ws.Cells.Item(10, 19).Value = 0.00191494505463948
$ws.Cells.Item(10, 20).Value = 0.00191494505463948

$ws.Cells.Item(11, 8).Value = 7.555711000000001
$ws.Cells.Item(11, 9).Value = 0.1726461783080517
$ws.Cells.Item(11, 10).Value = 0.1726461783080516
$ws.Cells.Item(11, 13).Value = 83.95738966666666
$ws.Cells.Item(11, 14).Value = 251.872169
$ws.Cells.Item(11, 15).Value = 0.386328906947963
$ws.Cells.Item(11, 16).Value = 0.386328906947963
$ws.Cells.Item(11, 17).Value = 211.4525908785732
$ws.Cells.Item(11, 18).Value = 1903.073317907159
$ws.Cells.Item(11, 19).Value = 0.06669820935449272
$ws.Cells.Item(11, 20).Value = 0.06669820935449271

$ws.Cells.Item(12, 7).Value = 4.235286666666666
$ws.Cells.Item(12, 8).Value = 12.70586
$ws.Cells.Item(12, 9).Value = 0.2903258437382188
$ws.Cells.Item(12, 10).Value = 0.2903258437382187
$ws.Cells.Item(12, 13).Value = 18.95273633333333
$ws.Cells.Item(12, 14).Value = 56.858209
$ws.Cells.Item(12, 15).Value = 0.08721078561875104
$ws.Cells.Item(12, 16).Value = 0.08721078561875105
$ws.Cells.Item(12, 17).Value = 80.27027148941553
$ws.Cells.Item(12, 18).Value = 722.43244340474
$ws.Cells.Item(12, 19).Value = 0.02531954491783681
$ws.Cells.Item(12, 20).Value = 0.02531954491783681

$ws.Cells.Item(13, 7).Value = 4.235286666666666
$ws.Cells.Item(13, 8).Value = 12.70586
$ws.Cells.Item(13, 9).Value = 0.2903258437382188
$ws.Cells.Item(13, 10).Value = 0.2903258437382187
$ws.Cells.Item(13, 15).Value = 0.04852204497892696
$ws.Cells.Item(13, 16).Value = 0.04852204497892696
$ws.Cells.Item(13, 17).Value = 44.66050495986666
$ws.Cells.Item(13, 18).Value = 401.9445446388
$ws.Cells.Item(13, 19).Value = 0.01408720364841077
$ws.Cells.Item(13, 20).Value = 0.01408720364841077

$ws.Cells.Item(14, 7).Value = 4.235286666666666
$ws.Cells.Item(14, 8).Value = 12.70586
$ws.Cells.Item(14, 9).Value = 0.2903258437382188
$ws.Cells.Item(14, 10).Value = 0.2903258437382187
$ws.Cells.Item(14, 13).Value = 101.4555613333333
$ws.Cells.Item(14, 14).Value = 304.366684
$ws.Cells.Item(14, 15).Value = 0.4668465309523581
$ws.Cells.Item(14, 16).Value = 0.4668465309523581
$ws.Cells.Item(14, 17).Value = 429.6933861742489
$ws.Cells.Item(14, 18).Value = 3867.24047556824
$ws.Cells.Item(14, 19).Value = 0.1355376129950039
$ws.Cells.Item(14, 20).Value = 0.1355376129950038

$ws.Cells.Item(15, 7).Value = 4.235286666666666
$ws.Cells.Item(15, 8).Value = 12.70586
$ws.Cells.Item(15, 9).Value = 0.2903258437382188
$ws.Cells.Item(15, 10).Value = 0.2903258437382187
$ws.Cells.Item(15, 13).Value = 2.410466333333333
$ws.Cells.Item(15, 14).Value = 7.231399000000001
$ws.Cells.Item(15, 15).Value = 0.01109173150200089
$ws.Cells.Item(15, 16).Value = 0.01109173150200089
$ws.Cells.Item(15, 17).Value = 10.20901592201555
$ws.Cells.Item(15, 18).Value = 91.88114329814
$ws.Cells.Item(15, 19).Value = 0.003220216306836191
$ws.Cells.Item(15, 20).Value = 0.00322021630683619

$ws.Cells.Item(16, 7).Value = 4.235286666666666
$ws.Cells.Item(16, 8).Value = 12.70586
$ws.Cells.Item(16, 9).Value = 0.2903258437382188
$ws.Cells.Item(16, 10).Value = 0.2903258437382187
$ws.Cells.Item(16, 13).Value = 83.95738966666666
$ws.Cells.Item(16, 14).Value = 251.872169
$ws.Cells.Item(16, 15).Value = 0.386328906947963
$ws.Cells.Item(16, 16).Value = 0.386328906947963
$ws.Cells.Item(16, 17).Value = 355.583613023371
$ws.Cells.Item(16, 18).Value = 3200.25251721034
$ws.Cells.Item(16, 19).Value = 0.1121612658701312
$ws.Cells.Item(16, 20).Value = 0.1121612658701312

$ws.Cells.Item(17, 7).Value = 1.937427333333333
$ws.Cells.Item(17, 8).Value = 5.812282
$ws.Cells.Item(17, 9).Value = 0.1328092451588843
$ws.Cells.Item(17, 10).Value = 0.1328092451588843
$ws.Cells.Item(17, 13).Value = 18.95273633333333
$ws.Cells.Item(17, 14).Value = 56.858209
$ws.Cells.Item(17, 15).Value = 0.08721078561875104
$ws.Cells.Item(17, 16).Value = 0.08721078561875105
$ws.Cells.Item(17, 17).Value = 36.71954941365977
$ws.Cells.Item(17, 18).Value = 330.475944722938
$ws.Cells.Item(17, 19).Value = 0.01158239860773961
$ws.Cells.Item(17, 20).Value = 0.01158239860773961

$ws.Cells.Item(18, 7).Value = 1.937427333333333
$ws.Cells.Item(18, 8).Value = 5.812282
$ws.Cells.Item(18, 9).Value = 0.1328092451588843
$ws.Cells.Item(18, 10).Value = 0.1328092451588843
$ws.Cells.Item(18, 15).Value = 0.04852204497892696
$ws.Cells.Item(18, 16).Value = 0.04852204497892696
$ws.Cells.Item(18, 17).Value = 20.42989999017333
$ws.Cells.Item(18, 18).Value = 183.86909991156
$ws.Cells.Item(18, 19).Value = 0.006444176167216722
$ws.Cells.Item(18, 20).Value = 0.006444176167216721

$ws.Cells.Item(19, 7).Value = 1.937427333333333
$ws.Cells.Item(19, 8).Value = 5.812282
$ws.Cells.Item(19, 9).Value = 0.1328092451588843
$ws.Cells.Item(19, 10).Value = 0.1328092451588843
$ws.Cells.Item(19, 13).Value = 101.4555613333333
$ws.Cells.Item(19, 14).Value = 304.366684
$ws.Cells.Item(19, 15).Value = 0.4668465309523581
$ws.Cells.Item(19, 16).Value = 0.4668465309523581
$ws.Cells.Item(19, 17).Value = 196.5627776458764
$ws.Cells.Item(19, 18).Value = 1769.064998812888
$ws.Cells.Item(19, 19).Value = 0.0620015353808264
$ws.Cells.Item(19, 20).Value = 0.06200153538082639

$ws.Cells.Item(20, 7).Value = 1.937427333333333
$ws.Cells.Item(20, 8).Value = 5.812282
$ws.Cells.Item(20, 9).Value = 0.1328092451588843
$ws.Cells.Item(20, 10).Value = 0.1328092451588843
$ws.Cells.Item(20, 13).Value = 2.410466333333333
$ws.Cells.Item(20, 14).Value = 7.231399000000001
$ws.Cells.Item(20, 15).Value = 0.01109173150200089
$ws.Cells.Item(20, 16).Value = 0.01109173150200089
$ws.Cells.Item(20, 17).Value = 4.670103360279778
$ws.Cells.Item(20, 18).Value = 42.030930242518
$ws.Cells.Item(20, 19).Value = 0.001473084488285757
$ws.Cells.Item(20, 20).Value = 0.001473084488285757

$ws.Cells.Item(21, 7).Value = 1.937427333333333
$ws.Cells.Item(21, 8).Value = 5.812282
$ws.Cells.Item(21, 9).Value = 0.1328092451588843
$ws.Cells.Item(21, 10).Value = 0.1328092451588843
$ws.Cells.Item(21, 13).Value = 83.95738966666666
$ws.Cells.Item(21, 14).Value = 251.872169
$ws.Cells.Item(21, 15).Value = 0.386328906947963
$ws.Cells.Item(21, 16).Value = 0.386328906947963
$ws.Cells.Item(21, 17).Value = 162.6613415755175
$ws.Cells.Item(21, 18).Value = 1463.952074179658
$ws.Cells.Item(21, 19).Value = 0.05130805051481582
$ws.Cells.Item(21, 20).Value = 0.05130805051481582

$ws.Cells.Item(22, 7).Value = 2.869646666666667
$ws.Cells.Item(22, 8).Value = 8.60894
$ws.Cells.Item(22, 9).Value = 0.1967122075319342
$ws.Cells.Item(22, 10).Value = 0.1967122075319342
$ws.Cells.Item(22, 13).Value = 18.95273633333333
$ws.Cells.Item(22, 14).Value = 56.858209
$ws.Cells.Item(22, 15).Value = 0.08721078561875104
$ws.Cells.Item(22, 16).Value = 0.08721078561875105
$ws.Cells.Item(22, 17).Value = 54.38765664316222
$ws.Cells.Item(22, 18).Value = 489.48890978846
$ws.Cells.Item(22, 19).Value = 0.01715542615965878
$ws.Cells.Item(22, 20).Value = 0.01715542615965878

$ws.Cells.Item(23, 7).Value = 2.869646666666667
$ws.Cells.Item(23, 8).Value = 8.60894
$ws.Cells.Item(23, 9).Value = 0.1967122075319342
$ws.Cells.Item(23, 10).Value = 0.1967122075319342
$ws.Cells.Item(23, 15).Value = 0.04852204497892696
$ws.Cells.Item(23, 16).Value = 0.04852204497892696
$ws.Cells.Item(23, 17).Value = 30.26002234946667
$ws.Cells.Item(23, 18).Value = 272.3402011452
$ws.Cells.Item(23, 19).Value = 0.009544878581768525
$ws.Cells.Item(23, 20).Value = 0.009544878581768525

$ws.Cells.Item(24, 7).Value = 2.869646666666667
$ws.Cells.Item(24, 8).Value = 8.60894
$ws.Cells.Item(24, 9).Value = 0.1967122075319342
$ws.Cells.Item(24, 10).Value = 0.1967122075319342
$ws.Cells.Item(24, 13).Value = 101.4555613333333
$ws.Cells.Item(24, 14).Value = 304.366684
$ws.Cells.Item(24, 15).Value = 0.4668465309523581
$ws.Cells.Item(24, 16).Value = 0.4668465309523581
$ws.Cells.Item(24, 17).Value = 291.1416133949955
$ws.Cells.Item(24, 18).Value = 2620.27452055496
$ws.Cells.Item(24, 19).Value = 0.09183441168226381
$ws.Cells.Item(24, 20).Value = 0.09183441168226381

$ws.Cells.Item(25, 7).Value = 2.869646666666667
$ws.Cells.Item(25, 8).Value = 8.60894
$ws.Cells.Item(25, 9).Value = 0.1967122075319342
$ws.Cells.Item(25, 10).Value = 0.1967122075319342
$ws.Cells.Item(25, 13).Value = 2.410466333333333
$ws.Cells.Item(25, 14).Value = 7.231399000000001
$ws.Cells.Item(25, 15).Value = 0.01109173150200089
$ws.Cells.Item(25, 16).Value = 0.01109173150200089
$ws.Cells.Item(25, 17).Value = 6.917186678562222
$ws.Cells.Item(25, 18).Value = 62.25468010706001
$ws.Cells.Item(25, 19).Value = 0.002181878989110092
$ws.Cells.Item(25, 20).Value = 0.002181878989110092

$ws.Cells.Item(26, 7).Value = 2.869646666666667
$ws.Cells.Item(26, 8).Value = 8.60894
$ws.Cells.Item(26, 9).Value = 0.1967122075319342
$ws.Cells.Item(26, 10).Value = 0.1967122075319342
$ws.Cells.Item(26, 13).Value = 83.95738966666666
$ws.Cells.Item(26, 14).Value = 251.872169
$ws.Cells.Item(26, 15).Value = 0.386328906947963
$ws.Cells.Item(26, 16).Value = 0.386328906947963
$ws.Cells.Item(26, 17).Value = 240.9280433989844
$ws.Cells.Item(26, 18).Value = 2168.35239059086
$ws.Cells.Item(26, 19).Value = 0.075995612119133
$ws.Cells.Item(26, 20).Value = 0.075995612119133
